$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.134.50'
$ws.Range("E2").Value = '  -0.10%  '
$ws.Range("D3").Value = '2.047.47'
$ws.Range("E3").Value = '  -1.42%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").Value = '''247.95'
$ws.Range("E5").Value = '  -2.27%  '
$ws.Range("D6").Value = '''0.662'
$ws.Range("E6").Value = '  -2.14%  '
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("D8").Value = '''56.10'
$ws.Range("E8").Value = '  -5.74%  '
$ws.Range("D9").Value = '''0.379'
$ws.Range("E9").Value = '  -3.35%  '
$ws.Range("D10").Value = '''0.0779'
$ws.Range("E10").Value = '  -2.78%  '
$ws.Range("D11").Value = '''0.108'
$ws.Range("E11").Value = '  -0.28%  '
$ws.Range("D12").Value = '''16.22'
$ws.Range("E12").Value = '  -0.78%  '
$ws.Range("D13").Value = '''0.878'
$ws.Range("E13").Value = '  +7.33%  '
$ws.Range("D14").Value = '2.351.49'
$ws.Range("E14").Value = '  -1.19%  '
$ws.Range("D15").Value = '''5.69'
$ws.Range("E15").Value = '  +2.02%  '
$ws.Range("D16").Value = '2.053.78'
$ws.Range("E16").Value = '  -1.12%  '
$ws.Range("D17").Value = '''18.36'
$ws.Range("E17").Value = '  +14.65%  '
$ws.Range("D18").Value = '37.143.34'
$ws.Range("E18").Value = '  -0.45%  '
$ws.Range("D19").Value = '''74.43'
$ws.Range("E19").Value = '  -0.59%  '
$ws.Range("D20").Value = '0.0₃0891'
$ws.Range("E20").Value = '  -4.38%  '
$ws.Range("D21").Value = '''5.38'
$ws.Range("E21").Value = '  -2.01%  '
$ws.Range("D22").Value = '''236.28'
$ws.Range("E22").Value = '  -1.35%  '
$ws.Range("D23").Value = '''0.999'
$ws.Range("E23").Value = '  -0.15%  '
$ws.Range("D24").Value = '''2.47'
$ws.Range("E24").Value = '  +2.41%  '
$ws.Range("D25").Value = '''9.50'
$ws.Range("E25").Value = '  +1.23%  '
$ws.Range("D26").Value = '''169.36'
$ws.Range("E26").Value = '  -0.57%  '
$ws.Range("E27").Value = '  -5.54%  '
$ws.Range("D28").Value = '''20.01'
$ws.Range("E28").Value = '  -2.10%  '
$ws.Range("E29").Value = '  -1.74%  '
$ws.Range("E30").Value = '  -1.25%  '
$ws.Range("D31").Value = '''4.83'
$ws.Range("E31").Value = '  +0.94%  '
$ws.Range("D32").Value = '''0.0617'
$ws.Range("E32").Value = '  -3.24%  '
$ws.Range("D33").Value = '''4.47'
$ws.Range("E33").Value = '  -0.95%  '
$ws.Range("D34").Value = '''0.0887'
$ws.Range("E34").Value = '  -2.98%  '
$ws.Range("E35").Value = '  +0.06%  '
$ws.Range("D36").Value = '''2.24'
$ws.Range("E36").Value = '  -3.36%  '
$ws.Range("D37").Value = '''1.77'
$ws.Range("E37").Value = '  -0.07%  '
$ws.Range("D38").Value = '''1.33'
$ws.Range("E38").Value = '  -3.10%  '
$ws.Range("D39").Value = '''5.26'
$ws.Range("E39").Value = '  +14.90%  '
$ws.Range("E40").Value = '  +9.15%  '
$ws.Range("D41").Value = '''0.0996'
$ws.Range("E41").Value = '  -15.40%  '
$ws.Range("D42").Value = '''0.0222'
$ws.Range("E42").Value = '  -2.41%  '
$ws.Range("E43").Value = '  -2.11%  '
$ws.Range("D44").Value = '''17.20'
$ws.Range("E44").Value = '  -4.27%  '
$ws.Range("D45").Value = '''95.32'
$ws.Range("E45").Value = '  -4.01%  '
$ws.Range("D46").Value = '''2.42'
$ws.Range("E46").Value = '  -3.63%  '
$ws.Range("D47").Value = '1.266.04'
$ws.Range("E47").Value = '  -3.42%  '
$ws.Range("E48").Value = '  -3.27%  '
$ws.Range("D49").Value = '''6.78'
$ws.Range("E49").Value = '  -2.53%  '
$ws.Range("D50").Value = '2.239.52'
$ws.Range("E50").Value = '  -1.09%  '
$ws.Range("D51").Value = '''43.81'
$ws.Range("E51").Value = '  -1.45%  '
